# Scraper run @ 31/12/2025 08:23:37 — adds the bus-arrival rows collected
# since the previous run (07:57:39) to the LP1912 sheet, mirrors the
# "215_EL PELIGRO" stop into the LP1912-215 sheet, and refreshes the
# "Última actualización" / "Total filas" header cells on every sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

$stamp = "Última actualización: 31/12/2025 08:23:37"

# --- Sheet "LP1912": 16 new scrape rows (744-759) -------------------------
# Columns: B=Hora_Scrap, C=Hora_Llegada, D=Línea, E=Minutos, F=Parada, G=Fecha
$newRows = @(
    @('08:23:26', '08:29', '14_ABASTO',       6,  'LP1912', '31/12/2025'),
    @('08:23:26', '08:44', '10_OLMOS',        21, 'LP1912', '31/12/2025'),
    @('08:23:26', '08:50', '23_HERNANDEZ',    27, 'LP1912', '31/12/2025'),
    @('08:23:26', '08:51', '16_SANTA ANA',    28, 'LP1912', '31/12/2025'),
    @('08:23:26', '09:02', '17X38_ROMERO',    39, 'LP1912', '31/12/2025'),
    @('08:23:26', '09:03', '23_HERNANDEZ',    40, 'LP1912', '31/12/2025'),
    @('08:23:26', '09:08', '16_SANTA ANA',    45, 'LP1912', '31/12/2025'),
    @('08:23:26', '09:14', '11_ETCHEVERRY',   51, 'LP1912', '31/12/2025'),
    @('08:23:26', '09:17', '27_EL RETIRO',    54, 'LP1912', '31/12/2025'),
    @('08:23:26', '09:21', '16_SANTA ANA',    58, 'LP1912', '31/12/2025'),
    @('08:23:26', '09:27', '215_EL PELIGRO',  64, 'LP1912', '31/12/2025'),
    @('08:23:26', '09:36', '23_HERNANDEZ',    73, 'LP1912', '31/12/2025'),
    @('08:23:26', '09:39', '15_ABASTO',       76, 'LP1912', '31/12/2025'),
    @('08:23:26', '09:44', '14_ABASTO',       81, 'LP1912', '31/12/2025'),
    @('08:23:26', '09:51', '15_ABASTO',       88, 'LP1912', '31/12/2025'),
    @('08:23:26', '09:58', '10_OLMOS',        95, 'LP1912', '31/12/2025')
)

$r = 744
foreach ($row in $newRows) {
    $ws1.Cells.Item($r, 2).Value = $row[0]
    $ws1.Cells.Item($r, 3).Value = $row[1]
    $ws1.Cells.Item($r, 4).Value = $row[2]
    $ws1.Cells.Item($r, 5).Value = $row[3]
    $ws1.Cells.Item($r, 6).Value = $row[4]
    $ws1.Cells.Item($r, 7).Value = $row[5]
    $r++
}

$ws1.Range("A2").Value = $stamp
$ws1.Range("A3").Value = "Total filas: 758"

# --- Sheet "LP1912-215": mirrors the 215_EL PELIGRO arrival as row 55 -----
# Columns: B=Fecha, C=Hora_Scrap, D=Hora_Llegada, E=Línea, F=Minutos, G=Parada
$ws2.Cells.Item(55, 2).Value = '31/12/2025'
$ws2.Cells.Item(55, 3).Value = '08:23:26'
$ws2.Cells.Item(55, 4).Value = '09:27'
$ws2.Cells.Item(55, 5).Value = '215_EL PELIGRO'
$ws2.Cells.Item(55, 6).Value = 64
$ws2.Cells.Item(55, 7).Value = 'LP1912'

$ws2.Range("A2").Value = $stamp
$ws2.Range("A3").Value = "Total filas: 54"

# --- Sheet "6203-6173": only the run timestamp refreshes, no new rows -----
$ws3.Range("A2").Value = $stamp
